# PPDM_Cookies Categories NameTest Cases_dataNupur.xlsx
# Rebuild the "Vendor name Data" sheet (sheet 2) as a two-column
# (Vendor name / Vendor Url) lookup table with 4 vendors, replacing the
# previous single "Vendor name data" column (which mixed names/urls).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Vendor"
$ws2 = $wb.Worksheets.Item(2)   # "Vendor name Data"

# --- Clear the old region (A1:C4) so none of the previous content lingers ---
$ws2.Range("A1:C6").Clear() | Out-Null

# --- Populate values. Order matters: it controls the order new shared
#     strings are appended in, which keeps the produced workbook as close
#     as possible to the reference edit. ---
$ws2.Range("C5").Value = "https://www.kellton.com/"
$ws2.Range("C6").Value = "https://mirus-it.co.uk/"

$ws2.Range("B3").Value = "One Trust"
$ws2.Range("B4").Value = "Bright Bridge"
$ws2.Range("B5").Value = "Kellton Tech"
$ws2.Range("B6").Value = "Mirus IT/MiContent Cloud"

$ws2.Range("B1").Value = "Vendor name"
$ws2.Range("C1").Value = "Vendor Url"

$ws2.Range("C3").Value = "https://www.onetrust.com/"
$ws2.Range("C4").Value = "https://brightbridgesolutions.com/"

$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("A5").Value = 3
$ws2.Range("A6").Value = 4

# --- Formatting ---

# Plain thin border (all sides), no fill -> used for A1 and the blank
# spacer row (row 2).
$plainBorder = $ws1.Range("B4")
$plainBorder.Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2:C2").PasteSpecial(-4122) | Out-Null

# Bold font + thin border (all sides), no fill -> header row labels.
$boldBorder = $ws1.Range("K2")
$boldBorder.Copy() | Out-Null
$ws2.Range("B1:C1").PasteSpecial(-4122) | Out-Null

# Thin border + centered/top-aligned -> row index column (A3:A6).
$plainBorder.Copy() | Out-Null
$ws2.Range("A3:A6").PasteSpecial(-4122) | Out-Null
$ws2.Range("A3:A6").HorizontalAlignment = -4108
$ws2.Range("A3:A6").VerticalAlignment = -4160

# Thin border + centered (horizontal & vertical) -> vendor name/url cells.
$plainBorder.Copy() | Out-Null
$ws2.Range("B3:C6").PasteSpecial(-4122) | Out-Null
$ws2.Range("B3:C6").HorizontalAlignment = -4108
$ws2.Range("B3:C6").VerticalAlignment = -4108

# Column C should be as wide as column B (both hold long text / URLs).
$ws2.Columns("C").ColumnWidth = $ws2.Columns("B").ColumnWidth

# Restore the selection/active cell bookkeeping to match the edited file:
# cursor left on C14 in the "Vendor name Data" sheet, but "Vendor" stays
# the active/visible tab.
$ws2.Range("C14").Select() | Out-Null
$ws1.Activate() | Out-Null
